# Generate Report for Handoff
#
# Updates the localization-status report to reflect a new handoff run:
#   - the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamps for the rows that were just (re-)handed off move forward
#     a few seconds
#   - those same rows now carry the "ht" (handoff-type) Priority marker
#     on the zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 14)

# Overview sheet: column G ("Latest HO Xliff Generate Date")
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-03 00:23:59"
}

# zh-cn sheet: column H ("Latest Handoff Datetime") + column E ("Priority")
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-09-03 00:23:54"
    $zhcn.Range("E$r").Value = "ht"
}

# de-de sheet: column H ("Latest Handoff Datetime") + column E ("Priority")
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-09-03 00:23:59"
    $dede.Range("E$r").Value = "ht"
}

